$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$aw = $excel.ActiveWindow

# --- Workbook view: tab ratio (bookViews/workbookView@tabRatio 655 -> 881) ---
$aw.TabRatio = 0.881

# --- Input cells: x, begin span, end span ---
$ws.Range("C2").Value = 75
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 100

# --- "cur axle location" column (C8:C25), each shifted down by 50 ---
$ws.Range("C8").Value = 131
$ws.Range("C9").Value = 123
$ws.Range("C10").Value = 118
$ws.Range("C11").Value = 113
$ws.Range("C12").Value = 108
$ws.Range("C13").Value = 99
$ws.Range("C14").Value = 94
$ws.Range("C15").Value = 88
$ws.Range("C16").Value = 83
$ws.Range("C17").Value = 75
$ws.Range("C18").Value = 67
$ws.Range("C19").Value = 62
$ws.Range("C20").Value = 57
$ws.Range("C21").Value = 52
$ws.Range("C22").Value = 43
$ws.Range("C23").Value = 38
$ws.Range("C24").Value = 32
$ws.Range("C25").Value = 27

$excel.Calculate()

# --- Sheet view: scroll position + selection ---
$aw.ScrollRow = 19
$aw.ScrollColumn = 1
$ws.Range("C3").Select()
